$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.214.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "'1.860.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'236.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.2872"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("D9").Value = "'0.06548"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "'21.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.01%  "
$ws.Range("D11").Value = "'0.07927"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "'97.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "'1.868.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "'0.6821"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "'268.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.91%  "
$ws.Range("D17").Value = "'30.218.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "'13.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.61%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'0.000007429"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("D21").Value = "'2.109.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "'5.329"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.57%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'6.187"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "'167.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").Value = "'9.228"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").Value = "'18.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "'1.964"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.53%  "
$ws.Range("D29").Value = "'1.384"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.43%  "
$ws.Range("D30").Value = "'0.09849"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("D31").Value = "'4.391"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").Value = "'1.475"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "'4.075"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").Value = "'0.04713"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "'1.135"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.57%  "
$ws.Range("D36").Value = "'0.7040"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").Value = "'2.705"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").Value = "'2.629"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.29%  "
$ws.Range("D40").Value = "'6.252"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("D41").Value = "'74.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").Value = "'1.943"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'0.8461"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "'0.9991"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "'103.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").Value = "'957.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("D48").Value = "'7.172"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "'9.236"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").Value = "'34.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "'0.05659"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.49%  "
